$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 573
    $ws.Range("F6").Value = 115
    $ws.Range("F10").Value = 4984
    $ws.Range("F11").Value = 4683
}
